# Generate Report for Handback
# Adds a new handback record (d5281af4-5f82-4ee8-adf9-65047c3bd308) as row 4
# to the Overview / zh-cn / de-de sheets (and their backing tables).

$wb = $excel.ActiveWorkbook

$fileId   = "d5281af4-5f82-4ee8-adf9-65047c3bd308"
$mdName   = "$fileId.md"
$mdPath   = "e2e\$fileId.md"
$inSync   = "Handed back: in sync with en-US"
$dotMd    = ".md"

$zhXlf    = "$fileId.11e66e2f369de819d7d1a081d2f8a271ed85448b.zh-cn.xlf"
$deXlf    = "$fileId.11e66e2f369de819d7d1a081d2f8a271ed85448b.de-de.xlf"

$genDate       = "2016-08-24 00:43:10"
$zhHoDate      = "2016-08-24 00:43:05"
$zhHbDate      = "2016-08-24 00:43:28"
$deHbDate      = "2016-08-24 00:43:36"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $mdPath
$wsOverview.Range("C4").Value = $dotMd
$wsOverview.Range("E4").Value = $inSync
$wsOverview.Range("F4").Value = $inSync
$wsOverview.Range("G4").NumberFormat = $dateFmt
$wsOverview.Range("G4").Value = $genDate

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e07e1dc4e5b6a1b5e6b1f9fd4e43b5fbd3c5a6e1/e2e/$mdName",
    "",
    "",
    $mdPath
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = $mdName
$wsZhCn.Range("B4").Value = $dotMd
$wsZhCn.Range("C4").Value = $inSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").NumberFormat = $dateFmt
$wsZhCn.Range("H4").Value = $zhHoDate
$wsZhCn.Range("I4").Value = $mdName
$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").NumberFormat = $dateFmt
$wsZhCn.Range("K4").Value = $zhHbDate
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e07e1dc4e5b6a1b5e6b1f9fd4e43b5fbd3c5a6e1/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c1a5a8b9e6f3d2c4b7a9e8f1d3c5b7a9e8f1d3c5/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = $mdName
$wsDeDe.Range("B4").Value = $dotMd
$wsDeDe.Range("C4").Value = $inSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").NumberFormat = $dateFmt
$wsDeDe.Range("H4").Value = $genDate
$wsDeDe.Range("I4").Value = $mdName
$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").NumberFormat = $dateFmt
$wsDeDe.Range("K4").Value = $deHbDate
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e07e1dc4e5b6a1b5e6b1f9fd4e43b5fbd3c5a6e1/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f8e7d6c5b4a3f2e1d0c9b8a7f6e5d4c3b2a1f0e/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null
